$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resize the affected ListObjects (Excel tables) to their new ranges
# first, while their current header text is still intact, so the table's
# column name is preserved across the move ---
$ws.ListObjects.Item("Tabla6").Resize($ws.Range("D2:D7"))
$ws.ListObjects.Item("Tabla13").Resize($ws.Range("D10:D14"))
$ws.ListObjects.Item("Tabla16").Resize($ws.Range("D16:D20"))
$ws.ListObjects.Item("Tabla18").Resize($ws.Range("D23:D27"))
$ws.ListObjects.Item("Tabla9").Resize($ws.Range("H11:H19"))

# --- Fix the "Datos" table (Tabla6): drop the stray "id_estacionamiento" row
# and shift "id_cordenadas" up, shrinking the table from D2:D8 to D2:D7 ---
$ws.Range("D7").Value2 = "id_cordenadas"
$ws.Range("D8").ClearContents()

# --- Move the "Servicios" table (Tabla13) up one row: D11:D15 -> D10:D14 ---
$ws.Range("D10").Value2 = "Servicios"
$ws.Range("D11").Value2 = "id_servicios"
$ws.Range("D12").Value2 = "nombre"
$ws.Range("D13").Value2 = "descripcion"
$ws.Range("D14").Value2 = "precio"
$ws.Range("D15").ClearContents()

# --- Extend the "estacionamiento" table (Tabla16) up one row and add a new
# "id_datos" column entry: D17:D20 -> D16:D20 ---
$ws.Range("D16").Value2 = "estacionamiento"
$ws.Range("D17").Value2 = "id_estacionamiento"
$ws.Range("D18").Value2 = "id_datos_estacionamiento"
$ws.Range("D19").Value2 = "id_datos"
$ws.Range("D20").Value2 = "id_servicios"

# --- Move the "Cordenadas" table (Tabla18) down one row: D22:D26 -> D23:D27 ---
$ws.Range("D22").ClearContents()
$ws.Range("D23").Value2 = "Cordenadas"
$ws.Range("D24").Value2 = "id_Cordenadas"
$ws.Range("D25").Value2 = "Cordenada X"
$ws.Range("D26").Value2 = "Cordenada Y"
$ws.Range("D27").Value2 = "id_tipo_cordenadas"

# --- Relocate the "Datos Estacionamiento" table (Tabla9) from F2:F10 to
# H11:H19, clearing out the old column F cells ---
$ws.Range("F2").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("F10").ClearContents()

$ws.Range("H11").Value2 = "Datos Estacionamiento"
$ws.Range("H12").Value2 = "id_datos_estacionamiento"
$ws.Range("H13").Value2 = "Nombre_estacionamiento"
$ws.Range("H14").Value2 = "id_ubicacion"
$ws.Range("H15").Value2 = "horario"
$ws.Range("H16").Value2 = "tarifas"
$ws.Range("H17").Value2 = "altura maxima"
$ws.Range("H18").Value2 = "descripcion"
$ws.Range("H19").Value2 = "id_asset (imagen representativa)"

# --- Column F is no longer auto-fit to the relocated table text, so widen
# it to the new fixed width ---
$ws.Columns("F").ColumnWidth = 26.85546875

# --- Update the active selection left behind by the edit ---
$ws.Range("D20").Select()
